$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Regenerated s_val data filtering out save games.
# Row 2
$ws.Range("B2").Value = [double]"0.1190320826869504"
$ws.Range("C2").Value = [double]"6.240767534437808e-05"
$ws.Range("D2").Value = [double]"0.7527432677738641"
$ws.Range("E2").Value = [double]"10.19245300693656"
$ws.Range("G2").Value = [double]"11.06429076507272"

# Row 3
$ws.Range("B3").Value = [double]"0.6606524410359556"
$ws.Range("C3").Value = [double]"1.655778082260271"
$ws.Range("D3").Value = [double]"3.537761648806719"
$ws.Range("E3").Value = [double]"10.19245300693656"
$ws.Range("G3").Value = [double]"16.0466451790395"

# Row 4
$ws.Range("B4").Value = [double]"1.455362044514542"
$ws.Range("C4").Value = [double]"1.655778082260271"
$ws.Range("D4").Value = [double]"0.1494219747398047"
$ws.Range("E4").Value = [double]"0.4942365360607697"
$ws.Range("G4").Value = [double]"3.754798637575387"

# Row 5
$ws.Range("B5").Value = [double]"0.1190320826869504"
$ws.Range("C5").Value = [double]"0.306821227259698"
$ws.Range("D5").Value = [double]"261.3203778131603"
$ws.Range("E5").Value = [double]"1133.036916526867"
$ws.Range("G5").Value = [double]"1394.783147649974"

# Row 6
$ws.Range("B6").Value = [double]"0.04271373187048222"
$ws.Range("C6").Value = [double]"0.306821227259698"
$ws.Range("D6").Value = [double]"0.1494219747398047"
$ws.Range("E6").Value = [double]"0.4942365360607697"
$ws.Range("G6").Value = [double]"0.9931934699307545"

# Row 7
$ws.Range("B7").Value = [double]"0.6606524410359556"
$ws.Range("C7").Value = [double]"1.655778082260271"
$ws.Range("D7").Value = [double]"0.1494219747398047"
$ws.Range("E7").Value = [double]"0.4942365360607697"
$ws.Range("G7").Value = [double]"2.960089034096801"

# Row 8
$ws.Range("B8").Value = [double]"1.455362044514542"
$ws.Range("C8").Value = [double]"1.655778082260271"
$ws.Range("D8").Value = [double]"0.7527432677738641"
$ws.Range("E8").Value = [double]"0.4942365360607697"
$ws.Range("G8").Value = [double]"4.358119930609447"

# Row 9
$ws.Range("B9").Value = [double]"3.286832544864788"
$ws.Range("C9").Value = [double]"1.655778082260271"
$ws.Range("D9").Value = [double]"0.7527432677738641"
$ws.Range("E9").Value = [double]"0.4942365360607697"
$ws.Range("G9").Value = [double]"6.189590430959694"

# Row 10
$ws.Range("B10").Value = [double]"1.455362044514542"
$ws.Range("C10").Value = [double]"1.655778082260271"
$ws.Range("D10").Value = [double]"3.537761648806719"
$ws.Range("E10").Value = [double]"10.19245300693656"
$ws.Range("G10").Value = [double]"16.84135478251809"
